$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 186, pushing existing rows 186-189 down to 187-190
$ws.Rows.Item(186).Insert()

# Populate the new row 186 with the new weekly record
$ws.Cells.Item(186, 1).Value = 7
$ws.Cells.Item(186, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(186, 3).Value = "Ñuble"
$ws.Cells.Item(186, 4).Value = 44595
$ws.Cells.Item(186, 5).Value = 16
$ws.Cells.Item(186, 6).Value = 100112032
$ws.Cells.Item(186, 7).Value = "Zapallo italiano"
$ws.Cells.Item(186, 8).Value = "Sin especificar"
$ws.Cells.Item(186, 9).Value = "Primera"
$ws.Cells.Item(186, 10).Value = 80
$ws.Cells.Item(186, 11).Value = 6000
$ws.Cells.Item(186, 12).Value = 6500
$ws.Cells.Item(186, 13).Value = 6250
$ws.Cells.Item(186, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(186, 15).Value = "Región del Maule"
$ws.Cells.Item(186, 16).Value = 104
$ws.Cells.Item(186, 17).Value = 60
$ws.Cells.Item(186, 18).Value = "Hortaliza"
